$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.9999960544946203"
$ws.Range("E2").Value = [double]"0.9999960544946203"

$ws.Range("D3").Value = [double]"0.9999999975531375"
$ws.Range("E3").Value = [double]"0.9999999975531375"

$ws.Range("D4").Value = [double]"0.9999900415990234"
$ws.Range("E4").Value = [double]"0.9999900415990234"

$ws.Range("D5").Value = [double]"0.02055988863629653"
$ws.Range("E5").Value = [double]"0.02055988863629653"

$ws.Range("D6").Value = [double]"1.844144426982434E-10"
$ws.Range("E6").Value = [double]"1.844144426982434E-10"

$ws.Range("D7").Value = [double]"1.039775356012195E-05"
$ws.Range("E7").Value = [double]"0.9999896022464398"

$ws.Range("D8").Value = [double]"0.9999999989235941"
$ws.Range("E8").Value = [double]"1.076405853339679E-09"

$ws.Range("D9").Value = [double]"1.026191607291095E-05"
$ws.Range("E9").Value = [double]"0.9999897380839271"

$ws.Range("D11").Value = [double]"0.008343507701603926"
$ws.Range("E11").Value = [double]"0.991656492298396"
$ws.Range("F11").Value = [double]"7.155652046203613"
